$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row 56 data (Segment Tree / Binary Indexed Tree question)
# NB: order of first-use determines shared-string insertion order, so
# write B, G, E (matching Range Sum Query -> lintcode 307 -> SegmentTree/BIT).
$ws.Cells.Item(56, 1).Value = 55
$ws.Cells.Item(56, 2).Value = "Range Sum Query - Mutable"
$ws.Cells.Item(56, 3).Value = 5
$ws.Cells.Item(56, 7).Value = "lintcode 307"
$ws.Cells.Item(56, 5).Value = "SegmentTree/BIT"
$ws.Cells.Item(56, 6).Value = "medium"

# Match existing style/alignment used by the other rows in this block
$ws.Cells.Item(56, 1).HorizontalAlignment = -4108  # xlCenter
$ws.Cells.Item(56, 2).HorizontalAlignment = -4131  # xlLeft
$ws.Cells.Item(56, 3).HorizontalAlignment = -4131  # xlLeft
$ws.Cells.Item(56, 5).HorizontalAlignment = -4108  # xlCenter
$ws.Cells.Item(56, 6).HorizontalAlignment = -4108  # xlCenter
$ws.Cells.Item(56, 7).HorizontalAlignment = -4108  # xlCenter

# Column E needs to grow to fit the new longer "SegmentTree/BIT" value
$ws.Columns.Item(5).EntireColumn.AutoFit()

# Update the active selection, matching the post-edit workbook state
$ws.Range("E61").Select()
